$d = $word.ActiveDocument

# Move the active window into the primary (default) header story and grab
# the HeaderFooter off the Selection -- going through Selection.HeaderFooter
# (as the Word UI's Insert > Header does) adds only the single "default"
# header part for this section instead of eagerly materializing the
# even-page / first-page header and footer siblings too.
$word.ActiveWindow.View.SeekView = 9   # wdSeekPrimaryHeader

$headerFooter = $word.Selection.HeaderFooter
$headerRange = $headerFooter.Range
$headerRange.Text = "Questionnaire 42"

# Paragraph-level formatting: the "Header" style, centered.
$headerRange.Paragraphs(1).Style = $d.Styles("Header")
$headerRange.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter

# Character-level formatting on just the typed text (exclude the trailing
# paragraph mark so the font change lands on the run, not the pilcrow).
$textRange = $headerRange.Duplicate
[void]$textRange.MoveEnd(1, -1)   # wdCharacter
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12

# Return the view to the main document body.
$word.ActiveWindow.View.SeekView = 0   # wdSeekMainDocument
